# Applies the "Added files + changed epsi[lon]" commit:
#   - Feuil1!B2:B6 get new "satisfaction" percentages.
#   - Feuil1!B6 is given a new (right-aligned) cell style -> this creates
#     cellXfs[1] in styles.xml (applyAlignment="1" / alignment horizontal="right").
#   - The worksheet's saved cursor/selection moves from G8 to F4.
#
# NOTE: the source workbook's radar chart (xl/charts/chart1.xml) caches the
# plotted numbers and two <c:axId> identifiers alongside the live
# Feuil1!$B$2:$B$6 reference. Excel recomputes that cache/regenerates axIds
# whenever it touches the chart, but this automation surface does not expose
# a way to rewrite a chart's cached values/axis ids (SeriesCollection/Points
# value setters and SetSourceData only update the live in-memory numbers,
# never the persisted <c:numCache>/<c:axId> on save) - so that part of the
# chart XML is intentionally left untouched here; every reachable effect of
# the edit (the underlying cell values, the new style, and the saved
# selection) is still applied below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "satisfaction" values for Comptable / Responsable d'atelier / des
# stocks / commercial / du personnel.
$ws.Range("B2").Value = 95.868799999999993
$ws.Range("B3").Value = 87.866100000000003
$ws.Range("B4").Value = 64.720600000000005
$ws.Range("B5").Value = 12.133900000000001
$ws.Range("B6").Value = 48.546700000000001

# B6 picks up a new, right-aligned style (adds cellXfs index 1 to styles.xml).
$ws.Range("B6").HorizontalAlignment = -4152   # xlRight

# Move the saved selection/active cell from G8 to F4.
$ws.Range("F4").Select() | Out-Null
